$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = -38915.87
$ws.Range("C4").Value = -42885.77
$ws.Range("D4").Value = -38293.06
$ws.Range("E4").Value = -31954.07
$ws.Range("F4").Value = -28973.95
$ws.Range("G4").Value = -181022.72
$ws.Range("B5").Value = 57327.05
$ws.Range("C5").Value = 69089.91
$ws.Range("D5").Value = 46505.97
$ws.Range("E5").Value = 68042.99000000001
$ws.Range("F5").Value = 97647.47
$ws.Range("G5").Value = 338613.39
$ws.Range("B6").Value = -36353.05
$ws.Range("C6").Value = -30457.38
$ws.Range("D6").Value = -38819.41
$ws.Range("E6").Value = -38946.09
$ws.Range("F6").Value = -34430.07
$ws.Range("G6").Value = -179006
$ws.Range("B7").Value = 452802.19
$ws.Range("C7").Value = 439948.04
$ws.Range("D7").Value = 521104.54
$ws.Range("E7").Value = 400785.61
$ws.Range("F7").Value = 439761.29
$ws.Range("G7").Value = 2254401.67
$ws.Range("B8").Value = -11866.61
$ws.Range("C8").Value = -20344.57
$ws.Range("D8").Value = -15438.74
$ws.Range("E8").Value = -15008.15
$ws.Range("F8").Value = -18061.25
$ws.Range("G8").Value = -80719.32000000001
$ws.Range("B9").Value = -42160.42
$ws.Range("C9").Value = -40308.66
$ws.Range("D9").Value = -31167.05
$ws.Range("E9").Value = -37792.51
$ws.Range("F9").Value = -39080.4
$ws.Range("G9").Value = -190509.04
$ws.Range("B10").Value = -11891.18
$ws.Range("C10").Value = -20277.56
$ws.Range("D10").Value = -16132.08
$ws.Range("E10").Value = -18620.19
$ws.Range("F10").Value = -13311.51
$ws.Range("G10").Value = -80232.52
$ws.Range("B11").Value = -57936.9
$ws.Range("C11").Value = -44474.84
$ws.Range("D11").Value = -54028.1
$ws.Range("E11").Value = -73444.08
$ws.Range("F11").Value = -49862.93
$ws.Range("G11").Value = -279746.85
$ws.Range("B12").Value = -28564.63
$ws.Range("C12").Value = -29994.23
$ws.Range("D12").Value = -27644.43
$ws.Range("E12").Value = -29930.08
$ws.Range("F12").Value = -24965.85
$ws.Range("G12").Value = -141099.22
$ws.Range("B13").Value = 196582.08
$ws.Range("C13").Value = 188143.44
$ws.Range("D13").Value = 185434.74
$ws.Range("E13").Value = 205682.74
$ws.Range("F13").Value = 189920.85
$ws.Range("G13").Value = 965763.85
$ws.Range("B14").Value = -196582.08
$ws.Range("C14").Value = -188143.44
$ws.Range("D14").Value = -185434.74
$ws.Range("E14").Value = -205682.74
$ws.Range("F14").Value = -189920.85
$ws.Range("G14").Value = -965763.85
$ws.Range("B15").Value = -18135.67
$ws.Range("C15").Value = -16482.58
$ws.Range("D15").Value = -15629.56
$ws.Range("E15").Value = -19907.67
$ws.Range("F15").Value = -12846.96
$ws.Range("G15").Value = -83002.44
$ws.Range("B16").Value = 81716.10000000001
$ws.Range("C16").Value = 76255.21000000001
$ws.Range("D16").Value = 76187.89999999999
$ws.Range("E16").Value = 68580.99000000001
$ws.Range("F16").Value = 110584.34
$ws.Range("G16").Value = 413324.54
$ws.Range("B17").Value = 346021.01
$ws.Range("C17").Value = 340067.57
$ws.Range("D17").Value = 406645.98
$ws.Range("E17").Value = 271806.75
$ws.Range("F17").Value = 426460.18
$ws.Range("G17").Value = 1791001.49
